$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.653.13"
$ws.Range("E2").Value = "  +0.54%  "

$ws.Range("D3").Value = "1.608.83"
$ws.Range("E3").Value = "  +0.40%  "

$ws.Range("E4").Value = "  -0.27%  "

$ws.Range("D5").Value = "'212.32"
$ws.Range("E5").Value = "  -0.35%  "

$ws.Range("E6").Value = "  +0.09%  "

$ws.Range("D7").Value = "'0.996"
$ws.Range("E7").Value = "  -0.23%  "

$ws.Range("D8").Value = "'28.89"
$ws.Range("E8").Value = "  +7.48%  "

$ws.Range("D9").Value = "'0.259"
$ws.Range("E9").Value = "  +3.55%  "

$ws.Range("E10").Value = "  +1.26%  "

$ws.Range("E11").Value = "  -0.72%  "

$ws.Range("D12").Value = "1.838.09"
$ws.Range("E12").Value = "  +0.39%  "

$ws.Range("D13").Value = "1.608.01"
$ws.Range("E13").Value = "  +0.64%  "

$ws.Range("E14").Value = "  +5.47%  "

$ws.Range("D15").Value = "'3.84"
$ws.Range("E15").Value = "  +2.70%  "

$ws.Range("D16").Value = "29.681.26"
$ws.Range("E16").Value = "  +0.49%  "

$ws.Range("B17").Value = "Chainlink"
$ws.Range("C17").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D17").Value = "'8.57"
$ws.Range("E17").Value = "  +12.52%  "

$ws.Range("B18").Value = "Litecoin"
$ws.Range("C18").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D18").Value = "'64.55"
$ws.Range("E18").Value = "  +1.58%  "

$ws.Range("D19").Value = "'241.46"
$ws.Range("E19").Value = "  +0.64%  "

$ws.Range("D20").Value = "0.0₃0703"
$ws.Range("E20").Value = "  +1.35%  "

$ws.Range("E21").Value = "  -0.27%  "

$ws.Range("E22").Value = "  +1.11%  "

$ws.Range("D23").Value = "'9.58"
$ws.Range("E23").Value = "  +4.44%  "

$ws.Range("D24").Value = "'2.10"
$ws.Range("E24").Value = "  +1.72%  "

$ws.Range("D25").Value = "'156.65"
$ws.Range("E25").Value = "  +1.26%  "

$ws.Range("D26").Value = "'15.56"
$ws.Range("E26").Value = "  +1.72%  "

$ws.Range("E27").Value = "  +1.00%  "

$ws.Range("D28").Value = "'6.54"
$ws.Range("E28").Value = "  +2.69%  "

$ws.Range("E29").Value = "  -0.21%  "

$ws.Range("E31").Value = "  +0.47%  "

$ws.Range("E32").Value = "  +0.91%  "

$ws.Range("E33").Value = "  +2.53%  "

$ws.Range("D34").Value = "1.425.17"
$ws.Range("E34").Value = "  +0.01%  "

$ws.Range("E35").Value = "  +5.11%  "

$ws.Range("E36").Value = "  +1.21%  "

$ws.Range("E37").Value = "  +2.32%  "

$ws.Range("E38").Value = "  -0.60%  "

$ws.Range("E39").Value = "  +2.84%  "

$ws.Range("D40").Value = "'0.552"
$ws.Range("E40").Value = "  +3.32%  "

$ws.Range("E41").Value = "  +5.51%  "

$ws.Range("E42").Value = "  +0.84%  "

$ws.Range("E43").Value = "  +3.86%  "

$ws.Range("D44").Value = "'54.27"
$ws.Range("E44").Value = "  +2.01%  "

$ws.Range("B45").Value = "Aave"
$ws.Range("C45").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D45").Value = "'68.34"
$ws.Range("E45").Value = "  +4.44%  "

$ws.Range("B46").Value = "PaxDollar"
$ws.Range("C46").Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$ws.Range("D46").Value = "'0.996"
$ws.Range("E46").Value = "  -0.28%  "

$ws.Range("D47").Value = "'0.998"
$ws.Range("E47").Value = "  +19.18%  "

$ws.Range("D48").Value = "'5.43"
$ws.Range("E48").Value = "  +2.62%  "

$ws.Range("D49").Value = "1.747.81"
$ws.Range("E49").Value = "  +0.31%  "

$ws.Range("D50").Value = "'87.30"
$ws.Range("E50").Value = "  +0.91%  "

$ws.Range("E51").Value = "  -1.33%  "
